# edit.ps1
# Updates the cryptocurrency list data (ranks 1-50) with refreshed market data
# as of 2023-10-01. Ranking is re-sorted by Market Cap (column E) descending,
# which shuffles several tickers/names between rows while updating
# Price, Market Cap, Volume and Change (24h) for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "BTC", "Bitcoin", 27148, 529466684050, 7335447880, 0.73496),
    @(3, "ETH", "Ethereum", 1685.26, 202673538772, 5919445066, 0.46486),
    @(4, "USDT", "Tether", 1, 83226160343, 13219102802, 0.01029),
    @(5, "BNB", "BNB", 215.88, 33219783537, 250734242, 0.51341),
    @(6, "XRP", "XRP", 0.519901, 27728182731, 482646947, 0.49799),
    @(7, "USDC", "USDC", 0.999266, 25098663953, 3454445392, -0.08739),
    @(8, "STETH", "Lido Staked Ether", 1685.45, 14818585841, 11242663, 0.53595),
    @(9, "SOL", "Solana", 22.99, 9497678799, 626214031, 13.01524),
    @(10, "ADA", "Cardano", 0.261157, 9171376999, 162181179, 4.1124),
    @(11, "DOGE", "Dogecoin", 0.062661, 8852289807, 201524618, 1.31411),
    @(12, "TRX", "TRON", 0.088904, 7915877087, 184908066, 0.17106),
    @(13, "TON", "Toncoin", 2.1, 7227839322, 17147204, -2.55712),
    @(14, "DOT", "Polkadot", 4.19, 5367637334, 89058084, 2.52768),
    @(15, "MATIC", "Polygon", 0.552614, 5144445816, 238031391, 4.6195),
    @(16, "LTC", "Litecoin", 67.52, 4970462410, 310908092, 2.77152),
    @(17, "BCH", "Bitcoin Cash", 236.71, 4614459960, 178071022, 0.37469),
    @(18, "LINK", "Chainlink", 8.11, 4517134590, 250192533, -0.95086),
    @(19, "WBTC", "Wrapped Bitcoin", 27202, 4429924911, 68158618, 0.7868),
    @(20, "SHIB", "Shiba Inu", 0.00000744, 4385265787, 82935896, 1.53487),
    @(21, "DAI", "Dai", 0.999956, 3811411700, 81941999, 0.0391),
    @(22, "TUSD", "TrueUSD", 0.998552, 3443358645, 161418993, -0.02249),
    @(23, "UNI", "Uniswap", 4.55, 3436812641, 81775290, 2.49827),
    @(24, "AVAX", "Avalanche", 9.6, 3407204364, 120468649, 4.71103),
    @(25, "LEO", "LEO Token", 3.67, 3404107808, 367579, -0.16756),
    @(26, "XLM", "Stellar", 0.113649, 3153938778, 41914333, 1.15305),
    @(27, "XMR", "Monero", 147.92, 2686242122, 54791971, 1.07574),
    @(28, "OKB", "OKB", 43.29, 2598185668, 2182104, 0.52967),
    @(29, "ETC", "Ethereum Classic", 16.43, 2353488191, 143436672, 2.30659),
    @(30, "BUSD", "BUSD", 0.999973, 2259921762, 2535888906, -0.01782),
    @(31, "ATOM", "Cosmos Hub", 7.28, 2131431692, 73008444, 0.97719),
    @(32, "HBAR", "Hedera", 0.050017, 1674609003, 17145550, 0.56653),
    @(33, "FIL", "Filecoin", 3.39, 1537129891, 64284138, 2.44599),
    @(34, "LDO", "Lido DAO", 1.68, 1498303539, 28593893, -0.17632),
    @(35, "ICP", "Internet Computer", 3.23, 1441226669, 16354090, 1.84671),
    @(36, "MKR", "Maker", 1558.04, 1405063488, 76955416, 5.18035),
    @(37, "CRO", "Cronos", 0.051144, 1345138792, 5904226, 0.81387),
    @(38, "QNT", "Quant", 90.98, 1324262915, 11333737, 0.62295),
    @(39, "APT", "Aptos", 5.53, 1315625053, 67202150, 2.65963),
    @(40, "MNT", "Mantle", 0.396476, 1282459757, 25708811, -0.64992),
    @(41, "VET", "VeChain", 0.01738331, 1264706234, 28513262, -0.91248),
    @(42, "ARB", "Arbitrum", 0.948768, 1211002777, 149551333, 4.88995),
    @(43, "NEAR", "NEAR Protocol", 1.15, 1126482424, 38626621, 1.62656),
    @(44, "OP", "Optimism", 1.38, 1105457918, 61155858, 4.6135),
    @(45, "KAS", "Kaspa", 0.04999857, 1047100321, 11646597, 3.02333),
    @(46, "AAVE", "Aave", 69.14, 1009315532, 85141078, 2.73448),
    @(47, "RETH", "Rocket Pool ETH", 1832.77, 970210527, 2026564, 0.611),
    @(48, "GRT", "The Graph", 0.091276, 845367723, 23879727, 4.03684),
    @(49, "ALGO", "Algorand", 0.1046, 820814306, 24458656, 2.15688),
    @(50, "WBT", "WhiteBIT Coin", 5.16, 743787475, 7077158, -0.51532),
    @(51, "USDD", "USDD", 0.997862, 726949805, 19499517, -0.058)
)

foreach ($row in $data) {
    $r      = $row[0]
    $ticker = $row[1]
    $name   = $row[2]
    $price  = $row[3]
    $mcap   = $row[4]
    $vol    = $row[5]
    $chg    = $row[6]

    $ws.Cells.Item($r, 2).Value = $ticker
    $ws.Cells.Item($r, 3).Value = $name
    $ws.Cells.Item($r, 4).Value = $price
    $ws.Cells.Item($r, 5).Value = $mcap
    $ws.Cells.Item($r, 6).Value = $vol
    $ws.Cells.Item($r, 7).Value = $chg
}
